$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (hyphens -> spaces)
$ws.Name = "GDS Leave orders from Jan 2024"

# Rewrite the header row: A1 becomes the new "SL.No" column, and
# B1:N1 get their own distinct header labels instead of all sharing "Name".
$ws.Range("A1").Value = "SL.No"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "designation"
$ws.Range("D1").Value = "officeName"
$ws.Range("E1").Value = "from"
$ws.Range("F1").Value = "to"
$ws.Range("G1").Value = "days"
$ws.Range("I1").Value = "accountNo"
$ws.Range("J1").Value = "remarks"
$ws.Range("H1").Value = "substituteName"
$ws.Range("K1").Value = "leaveType"
$ws.Range("L1").Value = "postmanBeatNo"
$ws.Range("M1").Value = "reference"
$ws.Range("N1").Value = "sendToHoOn"

# Fill in the new serial-number column (A2:A283) with 1..282
for ($i = 2; $i -le 283; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}
